# Weekly refresh of Albahaca price records (Vega Monumental Concepcion).
# Each data row (2-12, except 6 which is unchanged) is updated to reflect
# a different sampling date by pulling the Fecha/Volumen/Precio*/Origen
# values from another row in the same weekly cycle.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: replace with the date/price data previously on row 7
$ws.Range("D2").Value = 44650
$ws.Range("J2").Value = 130
$ws.Range("K2").Value = 3000
$ws.Range("L2").Value = 3500
$ws.Range("M2").Value = 3308
$ws.Range("P2").Value = 551

# Row 3: replace with the date/price data previously on row 5
$ws.Range("D3").Value = 44643
$ws.Range("J3").Value = 90
$ws.Range("K3").Value = 2800
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = 2911
$ws.Range("P3").Value = 485

# Row 4: replace with the date/price data previously on row 11
$ws.Range("D4").Value = 44672
$ws.Range("K4").Value = 3000
$ws.Range("L4").Value = 3500
$ws.Range("M4").Value = 3286
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 548

# Row 5: replace with the date/price data previously on row 3
$ws.Range("D5").Value = 44671
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 3500
$ws.Range("L5").Value = 4000
$ws.Range("M5").Value = 3733
$ws.Range("P5").Value = 622

# Row 7: replace with the date/price data previously on row 8
$ws.Range("D7").Value = 44685
$ws.Range("J7").Value = 150
$ws.Range("M7").Value = 3267
$ws.Range("P7").Value = 544

# Row 8: replace with the date/price data previously on row 9
$ws.Range("D8").Value = 44658
$ws.Range("J8").Value = 180
$ws.Range("K8").Value = 2500
$ws.Range("L8").Value = 3000
$ws.Range("M8").Value = 2778
$ws.Range("P8").Value = 463

# Row 9: replace with the date/price data previously on row 10
$ws.Range("D9").Value = 44631
$ws.Range("J9").Value = 110
$ws.Range("K9").Value = 3000
$ws.Range("L9").Value = 3500
$ws.Range("M9").Value = 3273
$ws.Range("O9").Value = "Provincia de Chacabuco"
$ws.Range("P9").Value = 546

# Row 10: replace with the date/price data previously on row 4
$ws.Range("D10").Value = 44644
$ws.Range("J10").Value = 140
$ws.Range("K10").Value = 2500
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = 2786
$ws.Range("P10").Value = 464

# Row 11: replace with the date/price data previously on row 12
$ws.Range("D11").Value = 44637
$ws.Range("J11").Value = 170
$ws.Range("K11").Value = 2800
$ws.Range("L11").Value = 3000
$ws.Range("M11").Value = 2906
$ws.Range("P11").Value = 484

# Row 12: replace with the date/price data previously on row 2
$ws.Range("D12").Value = 44659
$ws.Range("J12").Value = 90
$ws.Range("K12").Value = 2500
$ws.Range("M12").Value = 2722
$ws.Range("P12").Value = 454

